$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Sheet1!A2/A3 labels shortened from "S1"/"B1" to "S"/"B"
$ws1.Range("A2").Value = "S"
$ws1.Range("A3").Value = "B"

# Sheet1 becomes the active/selected sheet (previously Sheet2 was active),
# with C10 as the selected cell (previously C1).
$ws1.Activate()
$ws1.Range("C10").Select()
